$d = $word.ActiveDocument

# Locate the "# running the celery workers ..." paragraph (the last
# paragraph with real content before the trailing blank paragraphs) and
# the first trailing blank paragraph right after it.
$count = $d.Paragraphs.Count
$celeryIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*worker -l info --pool=solo -E*") {
        $celeryIndex = $i
        break
    }
}
if ($celeryIndex -eq -1) {
    $celeryIndex = $count - 2
}

# --- 1. Add a bottom paragraph border to that paragraph ---
$celeryPara = $d.Paragraphs.Item($celeryIndex)

$bdr = $celeryPara.Borders.Item(-3)
$bdr.LineStyle = 1
$bdr.LineWidth = 3
$bdr.Color = -16777216
$celeryPara.Borders.DistanceFromBottom = 1

# --- 2. Replace the first trailing empty paragraph with the new "7 Middleware" content ---
$targetPara = $d.Paragraphs.Item($celeryIndex + 1)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<w:p ' + $ns + '><w:r><w:t xml:space="preserve">7 Middleware  : </w:t></w:r></w:p>' + `
       '<w:p ' + $ns + '><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr>' + `
       '<w:r><w:t xml:space="preserve">We adding the middleware for counting the how many </w:t></w:r>' + `
       '<w:proofErr w:type="gramStart"/>' + `
       '<w:r><w:t>time</w:t></w:r>' + `
       '<w:proofErr w:type="gramEnd"/>' + `
       '<w:r><w:t xml:space="preserve"> requested on that path we are counting  </w:t></w:r></w:p>'

[void]$targetPara.Range.InsertXML($xml)
